$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# The old last data row ("Voting LogReg, Adaboost tree,Xgboost" submission 5)
# is dropped entirely - delete the worksheet row (also shrinks the table).
$ws.Rows.Item(6).Delete()

# Grow the table by one column (F) to hold the new "Training time" metric.
$lo.Resize($ws.Range("A1:F5"))

# --- Header row: renamed / reordered columns ---
$ws.Range("A1").Value = "Sumbission number"
$ws.Range("B1").Value = "Model"
$ws.Range("C1").Value = "F1_score_on_train_labelled"
$ws.Range("D1").Value = "F1_score_on_test_labelled"
$ws.Range("E1").Value = "F1_score_on_test_unlabelled"
$ws.Range("F1").Value = "Training_time_seconds"

# --- Row 2 : LogReg Multivariate ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "LogReg_Multivariate"
$ws.Range("C2").Value = 0.76545070847150998
$ws.Range("D2").Value = 0.75543478260869501
$ws.Range("E2").Value = 0.75461500000000004
$ws.Range("F2").Value = 0.4

# --- Row 3 : Adaboost with Decision Trees ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Adaboost_w_Decision_Trees"
$ws.Range("C3").Value = 0.76986754966887405
$ws.Range("D3").Value = 0.74481514878268695
$ws.Range("E3").Value = 0.75526700000000002
$ws.Range("F3").Value = 17.899999999999999

# --- Row 4 : Xgboost ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Xgboost"
$ws.Range("C4").Value = 0.77084112149532702
$ws.Range("D4").Value = 0.74895397489539695
$ws.Range("E4").Value = 0.75858400000000004
$ws.Range("F4").Value = 1.9

# --- Row 5 : Voting LogReg, Adaboost w/ DT, Xgboost ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Voting_LogReg_Adaboost_w_DT1_Xgboost"
$ws.Range("C5").Value = 0.76851364559055702
$ws.Range("D5").Value = 0.74849578820697904
$ws.Range("E5").Value = 0.75862099999999999
$ws.Range("F5").Value = 22.6

# Give the new column the same number-cell formatting as its neighbours.
$ws.Range("E2:E5").Copy()
$ws.Range("F2:F5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths (widened to fit the new, longer header names) ---
$ws.Columns.Item(2).ColumnWidth = 36.166666666666664
$ws.Columns.Item(3).ColumnWidth = 25.333333333333332
$ws.Columns.Item(4).ColumnWidth = 29
$ws.Columns.Item(5).ColumnWidth = 26.333333333333332
$ws.Columns.Item(6).ColumnWidth = 21

# --- Selection, matching the final saved cursor position ---
$ws.Range("D8").Select() | Out-Null
